$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header cells (bold, centered, bordered - style used by H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new data columns I and J
$valuesI = @(7, 9, 9, 8, 7, 6, 9)
$valuesJ = @(7, 9, 9, 8, 7, 7, 9)

for ($i = 0; $i -lt 7; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $valuesI[$i]
    $ws.Cells.Item($row, 10).Value = $valuesJ[$i]
}
